$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (shifts D:K -> E:L), to make room for the new quarter
$ws.Columns("D:D").Insert()

# Copy formatting from column E (the old column D, now shifted) into the newly inserted column D
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)

# Populate the newly inserted column D with the new quarter data
$ws.Range("D7").Value = 43407
$ws.Range("D8").Value = 92000
$ws.Range("D9").Value = 64600
$ws.Range("D10").Value = 27400
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = "NA"
$ws.Range("D15").Value = 2300
$ws.Range("D17").Value = 105700
$ws.Range("D18").Value = -13700
$ws.Range("D20").Value = 0
$ws.Range("D21").Value = -11400
$ws.Range("D22").Value = 300
$ws.Range("D23").Value = -14000
$ws.Range("D24").Value = 100
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = -14100
$ws.Range("D27").Value = -14100
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("D33").Value = -14100
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = -14100
$ws.Range("D38").Value = 43407
$ws.Range("D41").Value = 4500
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 5700
$ws.Range("D44").Value = 131300
$ws.Range("D45").Value = 13300
$ws.Range("D46").Value = 154800
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 12200
$ws.Range("D49").Value = 60200
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 11900
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 239100
$ws.Range("D57").Value = 42300
$ws.Range("D58").Value = 27400
$ws.Range("D59").Value = 15100
$ws.Range("D60").Value = 84800
$ws.Range("D61").Value = 0
$ws.Range("D62").Value = 25900
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 110600
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 15400
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 128400
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43407
$ws.Range("D81").Value = -14100
$ws.Range("D83").Value = 2300
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = -20400
$ws.Range("D91").Value = -1100
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -900
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 21100
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = -200
